# Incorporate updated data from upstream processes through 2024
# Update the "Solar" facility counts for year 2022 (row 24) and year 2024 (row 26).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24 corresponds to year 2022, column E ("Solar"): 26 -> 27
$ws.Range("E24").Value = 27

# Row 26 corresponds to year 2024, column E ("Solar"): 37 -> 48
$ws.Range("E26").Value = 48

$wb.Save()
